# C5-PowerPoint.pptx edit script
# 1) Re-style the table on slide 6 with the new table style GUID.
# 2) Swap the colour palette that the deck's (single) slide master /
#    design currently exposes ("Integral") for the palette that used to
#    live on the notes-master theme ("Office Theme") - i.e. apply the
#    Office Theme's 12 theme colours to the active theme colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Table style id change (slide 6, the table shape) -------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{502629F2-3AEB-4068-8741-87B7918BEA75}", $false)
    }
}

# --- 2. Theme colour palette swap -------------------------------------
# Office Theme colours (dk1..folHlink), expressed as OLE BGR-packed
# integers (0xBBGGRR) so that $color.RGB = ... writes the correct
# <a:srgbClr val="RRGGBB"/> into the theme part.
$officeThemeColors = @(
    0x000000,  # 1  dk1      000000
    0xFFFFFF,  # 2  lt1      FFFFFF
    0x6A5444,  # 3  dk2      44546A
    0xE6E6E7,  # 4  lt2      E7E6E6
    0xD59B5B,  # 5  accent1  5B9BD5
    0x317DED,  # 6  accent2  ED7D31
    0xA5A5A5,  # 7  accent3  A5A5A5
    0x00C0FF,  # 8  accent4  FFC000
    0xC47244,  # 9  accent5  4472C4
    0x47AD70,  # 10 accent6  70AD47
    0xC16305,  # 11 hlink    0563C1
    0x724F95   # 12 folHlink 954F72
)

$targetSlide = $p.Slides.Item(1)
$themeColors = $targetSlide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}

Write-Output "done"
